$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell content (set in the order that matches the target shared-string order) ---
$ws.Range("C7").Value = "입력문"

$ws.Range("D6").Value = 'printf("")'
$ws.Range("D6").Font.Size = 8

$ws.Range("D7").Value = 'scanf("")'
$ws.Range("D7").Font.Size = 8

$ws.Range("C8").Value = "조건문" + [char]10 + "conditional"
$ws.Range("C9").Value = "반복문" + [char]10 + "loop"

$ws.Range("D8").Value = "if" + [char]10 + "switch"
$ws.Range("D8").Font.Size = 8
$ws.Range("D8").WrapText = $true

$ws.Range("D9").Value = 'while (condition) { // Statement(s) }' + [char]10 + 'do {  // Statement(s) } while (condition);'
$ws.Range("D9").Font.Size = 8
$ws.Range("D9").WrapText = $true

# --- Row heights ---
$ws.Rows.Item(5).RowHeight = 34
$ws.Rows.Item(8).RowHeight = 34
$ws.Rows.Item(9).RowHeight = 34.5

# --- Column widths (closest achievable value; engine snaps width to a 7px/char grid) ---
$ws.Columns.Item(1).ColumnWidth = 3.857142857142857
$ws.Columns.Item(3).ColumnWidth = 12.428571428571429
$ws.Columns.Item(4).ColumnWidth = 30.0
$ws.Range("E1:J1").EntireColumn.ColumnWidth = 17.428571428571427

# --- Selection ---
$ws.Range("E9").Select() | Out-Null
